$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates -------------------------------------------------
# Some of the new Price values look like plain numbers (e.g. "598.52"). The sheet
# stores Price/Volume cells as text, so for those cells we briefly switch the
# cell to Text format before writing the value (otherwise Excel would silently
# reinterpret the text as a number and, e.g., drop a trailing zero). The style is
# restored to Normal immediately afterwards so only the cell value/type changes.
$priceCellsLookingNumeric = @('D5', 'D6', 'D10', 'D13', 'D14', 'D19', 'D20', 'D21', 'D22', 'D25', 'D26', 'D30', 'D36', 'D38', 'D39', 'D42', 'D47', 'D51')
$priceValues = @{
    'D5' = '598.52'
    'D6' = '159.77'
    'D10' = '0.141'
    'D13' = '0.360'
    'D14' = '28.24'
    'D19' = '11.85'
    'D20' = '364.81'
    'D21' = '7.63'
    'D22' = '4.53'
    'D25' = '74.32'
    'D26' = '1.00'
    'D30' = '583.95'
    'D36' = '1.64'
    'D38' = '160.66'
    'D39' = '19.81'
    'D42' = '5.37'
    'D47' = '157.40'
    'D51' = '21.98'
}
foreach ($ref in $priceCellsLookingNumeric) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = '@'
    $cell.Value = $priceValues[$ref]
    $cell.Style = 'Normal'
}

# --- Remaining Price (D) and Volume(1h) (E) updates -----------------------------
# These values already remain text on assignment (non-numeric strings), so they
# can be written directly.
$ws.Range('D2').Value = '68.511.41'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '2.694.18'
$ws.Range('E3').Value = '  +1.88%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('E6').Value = '  +2.05%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').Value = '2.693.29'
$ws.Range('E9').Value = '  +1.89%  '
$ws.Range('E10').Value = '  +0.27%  '
$ws.Range('E11').Value = '  -0.34%  '
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('E13').Value = '  +2.63%  '
$ws.Range('E14').Value = '  +0.91%  '
$ws.Range('D15').Value = '3.186.38'
$ws.Range('E15').Value = '  +1.93%  '
$ws.Range('E16').Value = '  -0.83%  '
$ws.Range('D17').Value = '68.464.14'
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('D18').Value = '2.693.15'
$ws.Range('E18').Value = '  +1.59%  '
$ws.Range('E19').Value = '  +4.15%  '
$ws.Range('E20').Value = '  +0.45%  '
$ws.Range('E21').Value = '  +3.86%  '
$ws.Range('E22').Value = '  +2.59%  '
$ws.Range('E23').Value = '  +1.97%  '
$ws.Range('E24').Value = '  +1.75%  '
$ws.Range('E25').Value = '  -1.50%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +1.20%  '
$ws.Range('E29').Value = '  +0.27%  '
$ws.Range('E30').Value = '  +4.94%  '
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('E32').Value = '  +1.63%  '
$ws.Range('E33').Value = '  +2.26%  '
$ws.Range('E34').Value = '  +5.19%  '
$ws.Range('E35').Value = '  +3.55%  '
$ws.Range('E36').Value = '  +5.90%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('E38').Value = '  +0.12%  '
$ws.Range('E39').Value = '  +0.46%  '
$ws.Range('E40').Value = '  +1.85%  '
$ws.Range('E41').Value = '  +2.22%  '
$ws.Range('E42').Value = '  +0.83%  '
$ws.Range('E43').Value = '  +2.02%  '
$ws.Range('E44').Value = '  +0.29%  '
$ws.Range('E47').Value = '  -0.87%  '
$ws.Range('E48').Value = '  +5.67%  '
$ws.Range('E49').Value = '  +4.80%  '
$ws.Range('E50').Value = '  +6.89%  '
$ws.Range('E51').Value = '  +0.08%  '
